$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in column D stay as text by
# pre-formatting the target cells as Text before assigning their values.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.700.46'
$ws.Range('E2').Value = '  -0.05%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.894.28'
$ws.Range('E3').Value = '  +1.28%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.26%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.06'
$ws.Range('E5').Value = '  +0.26%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.15%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4577'
$ws.Range('E7').Value = '  -1.44%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3869'
$ws.Range('E8').Value = '  -1.21%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.67'
$ws.Range('E9').Value = '  +0.66%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07877'
$ws.Range('E10').Value = '  -0.32%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  +3.24%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.71'
$ws.Range('E12').Value = '  -2.59%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.887.24'
$ws.Range('E13').Value = '  +0.73%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.032'
$ws.Range('E14').Value = '  +1.36%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.699'
$ws.Range('E15').Value = '  -0.74%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06952'
$ws.Range('E16').Value = '  -0.43%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.23'
$ws.Range('E17').Value = '  -1.33%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  +0.31%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001001'
$ws.Range('E19').Value = '  -0.80%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.01'
$ws.Range('E20').Value = '  +0.53%  '

# Row 21
$ws.Range('E21').Value = '  +0.15%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.710.84'
$ws.Range('E22').Value = '  -0.01%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.305'
$ws.Range('E23').Value = '  -0.52%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.00'
$ws.Range('E24').Value = '  -1.03%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.133.14'
$ws.Range('E25').Value = '  +2.84%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.062'
$ws.Range('E26').Value = '  -2.34%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.68'
$ws.Range('E27').Value = '  +0.76%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.30'
$ws.Range('E28').Value = '  -0.40%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.806'
$ws.Range('E29').Value = '  +1.28%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.33'
$ws.Range('E30').Value = '  -0.95%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.907'
$ws.Range('E31').Value = '  -4.67%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09308'
$ws.Range('E32').Value = '  -0.68%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9216'
$ws.Range('E33').Value = '  -1.67%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.313'
$ws.Range('E34').Value = '  -0.20%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.331'
$ws.Range('E35').Value = '  -1.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.265'
$ws.Range('E36').Value = '  -2.76%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05734'
$ws.Range('E37').Value = '  -1.97%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.160'
$ws.Range('E38').Value = '  +1.11%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02060'
$ws.Range('E39').Value = '  -3.08%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.716'
$ws.Range('E40').Value = '  -2.28%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5617'
$ws.Range('E41').Value = '  -0.74%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1782'
$ws.Range('E42').Value = '  -0.30%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.730'
$ws.Range('E43').Value = '  -2.39%  '

# Row 44
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.07171'
$ws.Range('E44').Value = '  -0.99%  '

# Row 45
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.197'
$ws.Range('E45').Value = '  +2.64%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.69'
$ws.Range('E46').Value = '  -0.62%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5311'
$ws.Range('E47').Value = '  -0.09%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.121'
$ws.Range('E48').Value = '  -1.36%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.826'
$ws.Range('E49').Value = '  -1.30%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '112.34'
$ws.Range('E50').Value = '  -1.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.459'
$ws.Range('E51').Value = '  +4.59%  '
